$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 00:22"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 789234
$ws.Range("C4").Value = 24598
$ws.Range("D4").Value = 71832
$ws.Range("E4").Value = 675104
$ws.Range("F4").Value = 13634
$ws.Range("G4").Value = 1723
$ws.Range("H4").Value = 42298

# Row 14: Brasil
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").Value = 40581
$ws.Range("C14").Value = 1927
$ws.Range("D14").Value = 22130
$ws.Range("E14").Value = 15876
$ws.Range("F14").Value = 7919
$ws.Range("G14").Value = 113
$ws.Range("H14").Value = 2575

# Row 50: Colombia
$ws.Range("A50").Value = "Colombia"
$ws.Range("B50").Value = 3977
$ws.Range("C50").Value = 185
$ws.Range("D50").Value = 804
$ws.Range("E50").Value = 2984
$ws.Range("F50").Value = 98
$ws.Range("G50").Value = 10
$ws.Range("H50").Value = 189

# Row 51: Finlandia
$ws.Range("A51").Value = "Finlandia"
$ws.Range("B51").Value = 3868
$ws.Range("C51").Value = 85
$ws.Range("D51").Value = 2000
$ws.Range("E51").Value = 1770
$ws.Range("F51").Value = 67
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 98

# Row 123: Somalia
$ws.Range("A123").Value = "Somalia"
$ws.Range("B123").Value = 237
$ws.Range("C123").Value = 73
$ws.Range("D123").Value = 4
$ws.Range("E123").Value = 225
$ws.Range("F123").Value = 2
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 8

# Row 124: Mali
$ws.Range("A124").Value = "Mali"
$ws.Range("B124").Value = 224
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 42
$ws.Range("E124").Value = 168
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 14

# Row 125: El Salvador
$ws.Range("A125").Value = "El Salvador"
$ws.Range("B125").Value = 218
$ws.Range("C125").Value = 17
$ws.Range("D125").Value = 46
$ws.Range("E125").Value = 165
$ws.Range("F125").Value = 2
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 7

# Row 126: Paraguay
$ws.Range("A126").Value = "Paraguay"
$ws.Range("B126").Value = 208
$ws.Range("C126").Value = 2
$ws.Range("D126").Value = 46
$ws.Range("E126").Value = 154
$ws.Range("F126").Value = 1
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 8

# Row 127: Jamaica
$ws.Range("A127").Value = "Jamaica"
$ws.Range("B127").Value = 196
$ws.Range("C127").Value = 23
$ws.Range("D127").Value = 27
$ws.Range("E127").Value = 164
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 5

# Row 128: Islas Feroe
$ws.Range("A128").Value = "Islas Feroe"
$ws.Range("B128").Value = 185
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 176
$ws.Range("E128").Value = 9
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0

# Row 144: Monaco
$ws.Range("A144").Value = "Monaco"
$ws.Range("B144").Value = 94
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 23
$ws.Range("E144").Value = 68
$ws.Range("F144").Value = 3
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 3

# Row 158: Polinesia Francesa
$ws.Range("A158").Value = "Polinesia Francesa"
$ws.Range("B158").Value = 56
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 19
$ws.Range("E158").Value = 37
$ws.Range("F158").Value = 1
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0

# Row 159: Uganda
$ws.Range("A159").Value = "Uganda"
$ws.Range("B159").Value = 56
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 38
$ws.Range("E159").Value = 18
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0
